$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Cluster Name"
$ws.Range("B1").Value = "Active cases"

$ws.Range("A2").Value = '3321 Rochester and Elmore District Health Service Yalunkan Aged Care Hostel Rochester'
$ws.Range("B2").Value = 10
$ws.Range("A3").Value = '3323 Villa Maria Catholic Homes St Bernadette''s Aged Care Sunshine North'
$ws.Range("B3").Value = 14
$ws.Range("A4").Value = '3600 Belvedere Aged Care Noble Park'
$ws.Range("B4").Value = 22
$ws.Range("A5").Value = '3601 Baptcare Westhaven community'
$ws.Range("B5").Value = 24
$ws.Range("A6").Value = '3653 Fronditha Thalpori St Albans Aged'
$ws.Range("B6").Value = 23
$ws.Range("A7").Value = '3939 Bupa Aged Care Eastwood'
$ws.Range("B7").Value = 15
$ws.Range("A8").Value = '3988 Kerala Manor Aged Care Diamond Creek'
$ws.Range("B8").Value = 10
$ws.Range("A9").Value = '4257 BlueCross The Gables Camberwell'
$ws.Range("B9").Value = 27
$ws.Range("A10").Value = '44087 Fitzroy Primary School Fitzroy'
$ws.Range("B10").Value = 11
$ws.Range("A11").Value = '44098 Stawell Primary School'
$ws.Range("B11").Value = 22
$ws.Range("A12").Value = '44121 Wallan Primary School Wallan'
$ws.Range("B12").Value = 12
$ws.Range("A13").Value = '44165 Greenvale Primary School'
$ws.Range("B13").Value = 17
$ws.Range("A14").Value = '44234 Lucknow Primary School Bairnsdale'
$ws.Range("B14").Value = 18
$ws.Range("A15").Value = '44444 Nar Nar Goon Primary School Nar Nar Goon'
$ws.Range("B15").Value = 14
$ws.Range("A16").Value = '44630 Black Rock Primary School Black Rock'
$ws.Range("B16").Value = 10
$ws.Range("A17").Value = '44667 Beaumaris Primary School Beaumaris'
$ws.Range("B17").Value = 16
$ws.Range("A18").Value = '44718 Parkdale Primary School Parkdale'
$ws.Range("B18").Value = 12
$ws.Range("A19").Value = '44811 Dandenong North Primary School Dandenong'
$ws.Range("B19").Value = 16
$ws.Range("A20").Value = '44812 Bairnsdale West Primary School'
$ws.Range("B20").Value = 10
$ws.Range("A21").Value = '44865 Parktone Primary School Parkdale'
$ws.Range("B21").Value = 28
$ws.Range("A22").Value = '44891 Cranbourne Park Primary School Cranbourne'
$ws.Range("B22").Value = 11
$ws.Range("A23").Value = '45248 Brookside P-9 College Caroline Springs'
$ws.Range("B23").Value = 19
$ws.Range("A24").Value = '45249 Creekside K-9 College Caroline Springs'
$ws.Range("B24").Value = 10
$ws.Range("A25").Value = '45267 Epping Views Primary School'
$ws.Range("B25").Value = 17
$ws.Range("A26").Value = '45518 Ashwood High School Ashwood'
$ws.Range("B26").Value = 22
$ws.Range("A27").Value = '45569 Nhill College Nhill'
$ws.Range("B27").Value = 40
$ws.Range("A28").Value = '45648 St Brendans Primary School Shepparton'
$ws.Range("B28").Value = 11
$ws.Range("A29").Value = '4574 Village Glen Aged Care Residences Mornington'
$ws.Range("B29").Value = 14
$ws.Range("A30").Value = '45784 Holy Rosary Primary School White Hills'
$ws.Range("B30").Value = 48
$ws.Range("A31").Value = '46037 Nazareth Catholic Primary School Grovedale'
$ws.Range("B31").Value = 24
$ws.Range("A32").Value = '46050 Our Lady''s Catholic Primary School Craigieburn'
$ws.Range("B32").Value = 17
$ws.Range("A33").Value = '46093 St Brendan''s Primary School Somerville'
$ws.Range("B33").Value = 12
$ws.Range("A34").Value = '46095 Bethany Catholic Primary School Werribee'
$ws.Range("B34").Value = 16
$ws.Range("A35").Value = '46125 Our Lady of the Southern Cross Primary School Manor Lakes'
$ws.Range("B35").Value = 15
$ws.Range("A36").Value = '46276 Hillcrest Christian College Clyde North'
$ws.Range("B36").Value = 14
$ws.Range("A37").Value = '46328 Ilim College Inverloch Crescent Dallas'
$ws.Range("B37").Value = 17
$ws.Range("A38").Value = '46390 Al Siraat College Epping'
$ws.Range("B38").Value = 37
$ws.Range("A39").Value = '50584 St Mary of the Cross MacKillop Primary School Epping'
$ws.Range("B39").Value = 10
$ws.Range("A40").Value = '52380 Al Iman College Melton South'
$ws.Range("B40").Value = 11
$ws.Range("A41").Value = '52473 John Henry Primary School Pakenham'
$ws.Range("B41").Value = 14
$ws.Range("A42").Value = 'Adass Israel School Elsternwick'
$ws.Range("B42").Value = 10
$ws.Range("A43").Value = 'Camp Coolamatong Farm Camp Banksia Peninsula'
$ws.Range("B43").Value = 13
$ws.Range("A44").Value = 'Christ the Priest Primary School Caroline Springs'
$ws.Range("B44").Value = 13
$ws.Range("A45").Value = 'Escala NewQuay Construction Site Docklands Drive Docklands'
$ws.Range("B45").Value = 12
$ws.Range("A46").Value = 'Hamilton Country Music Festival Hamilton Golf Club Hamilton'
$ws.Range("B46").Value = 29
$ws.Range("A47").Value = 'Islamic College of Melbourne Tarneit Oct Nov'
$ws.Range("B47").Value = 13
$ws.Range("A48").Value = 'Oakleigh Grammar Melbourne Private School Oakleigh'
$ws.Range("B48").Value = 11
$ws.Range("A49").Value = 'Social Gathering 20 November Sunbury'
$ws.Range("B49").Value = 13
$ws.Range("A50").Value = 'St Josephs Catholic Primary School Warragul'
$ws.Range("B50").Value = 13
$ws.Range("A51").Value = 'Wagstaff Meat Processing Plant Cranbourne East'
$ws.Range("B51").Value = 34
$ws.Range("A52").Value = 'Werribee Mercy Hospital Emergency Department'
$ws.Range("B52").Value = 12
